# Preference Logic Bug Fixed
# Swap employee names that were incorrectly assigned due to a preference-logic bug.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# WED 10-5 shift: was "Minal", should be "Minjung"
$ws.Range("C6").Value = "Minjung"

# SAT 10-5 shift: was "Minjung", should be "Minal"
$ws.Range("C12").Value = "Minal"

# SAT 5-11:30 A shift: was "Seoyoon", should be "yujin"
$ws.Range("C14").Value = "yujin"

# SAT 5-11:30 B shift: was "yujin", should be "Seoyoon"
$ws.Range("C15").Value = "Seoyoon"
